$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Die unten stehenden Teams machen den Anfang.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Teams, welche Im Folgenden aufgelistet sind, haben den Weg geebnet und das Projekt ins Leben gerufen.",
    2
)

$d.Content.Find.Execute(
    "Dieses Team konzentriert sich auf Gemeinschaftsbildung, Wachstum und allgemeine Nutzerakquise.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dieses Team richtet seinen Fokus auf Communitybildung, Wachstum und dem Anwerben neuer Mitglieder.",
    2
)

$d.Content.Find.Execute(
    "Verantwortlich für die Erstellung von Web Applications.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Verantwortlich das Erstellen von Web-Applikationen.",
    2
)

$d.Content.Find.Execute(
    "Sicherstellung der Erfüllung von Qualitätskriterien aller Entwicklungsaufgaben.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Stellt sicher, dass alle Qualitätskriterien erfüllt werden.",
    2
)
